# Rename the "Collection_DM" worksheet tab to "CRF_DM".
#
# Renaming the sheet through the Excel object model keeps everything that
# depends on the sheet name in sync automatically - in particular the
# workbook-scoped defined name `_xlnm._FilterDatabase`, whose formula
# reference ("Collection_DM!$A$1:$AK$16") must follow the sheet's new name
# ("CRF_DM!$A$1:$AK$16").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_DM")
$ws.Name = "CRF_DM"
